{"js": "// Update the worksheet date and every \"A\u00d7B=C\" answer cell to the new\n// values from the commit. Each old value is unique in the document, so a\n// straightforward search-and-replace per pair is unambiguous.\nconst pairs = [\n  [\"2025-05-20 Tuesday\", \"2025-05-21 Wednesday\"],\n  [\"289\u00d73=867\", \"586\u00d77=4102\"],\n  [\"725\u00d78=5800\", \"420\u00d78=3360\"],\n  [\"612\u00d73=1836\", \"998\u00d77=6986\"],\n  [\"957\u00d76=5742\", \"791\u00d74=3164\"],\n  [\"528\u00d76=3168\", \"269\u00d73=807\"],\n  [\"257\u00d72=514\", \"127\u00d75=635\"],\n  [\"805\u00d79=7245\", \"636\u00d76=3816\"],\n  [\"519\u00d74=2076\", \"783\u00d77=5481\"],\n  [\"783\u00d75=3915\", \"432\u00d79=3888\"],\n  [\"378\u00d77=2646\", \"563\u00d73=1689\"],\n  [\"597\u00d78=4776\", \"665\u00d76=3990\"],\n  [\"196\u00d74=784\", \"891\u00d78=7128\"],\n  [\"635\u00d79=5715\", \"447\u00d72=894\"],\n  [\"896\u00d75=4480\", \"346\u00d76=2076\"],\n  [\"948\u00d75=4740\", \"759\u00d78=6072\"],\n  [\"435\u00d76=2610\", \"908\u00d77=6356\"],\n  [\"941\u00d75=4705\", \"898\u00d77=6286\"],\n  [\"185\u00d78=1480\", \"675\u00d79=6075\"],\n  [\"553\u00d77=3871\", \"923\u00d72=1846\"],\n  [\"903\u00d76=5418\", \"555\u00d79=4995\"],\n  [\"105\u00d73=315\", \"720\u00d77=5040\"],\n  [\"831\u00d72=1662\", \"478\u00d75=2390\"],\n  [\"199\u00d74=796\", \"840\u00d76=5040\"],\n  [\"849\u00d74=3396\", \"575\u00d73=1725\"],\n  [\"795\u00d75=3975\", \"619\u00d79=5571\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"A\u00d7B=C\" answer cell to the new\n# values from the commit. Each old value is unique in the document, so a\n# straightforward Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2025-05-20 Tuesday\", \"2025-05-21 Wednesday\")\n    ,@(\"289\u00d73=867\", \"586\u00d77=4102\")\n    ,@(\"725\u00d78=5800\", \"420\u00d78=3360\")\n    ,@(\"612\u00d73=1836\", \"998\u00d77=6986\")\n    ,@(\"957\u00d76=5742\", \"791\u00d74=3164\")\n    ,@(\"528\u00d76=3168\", \"269\u00d73=807\")\n    ,@(\"257\u00d72=514\", \"127\u00d75=635\")\n    ,@(\"805\u00d79=7245\", \"636\u00d76=3816\")\n    ,@(\"519\u00d74=2076\", \"783\u00d77=5481\")\n    ,@(\"783\u00d75=3915\", \"432\u00d79=3888\")\n    ,@(\"378\u00d77=2646\", \"563\u00d73=1689\")\n    ,@(\"597\u00d78=4776\", \"665\u00d76=3990\")\n    ,@(\"196\u00d74=784\", \"891\u00d78=7128\")\n    ,@(\"635\u00d79=5715\", \"447\u00d72=894\")\n    ,@(\"896\u00d75=4480\", \"346\u00d76=2076\")\n    ,@(\"948\u00d75=4740\", \"759\u00d78=6072\")\n    ,@(\"435\u00d76=2610\", \"908\u00d77=6356\")\n    ,@(\"941\u00d75=4705\", \"898\u00d77=6286\")\n    ,@(\"185\u00d78=1480\", \"675\u00d79=6075\")\n    ,@(\"553\u00d77=3871\", \"923\u00d72=1846\")\n    ,@(\"903\u00d76=5418\", \"555\u00d79=4995\")\n    ,@(\"105\u00d73=315\", \"720\u00d77=5040\")\n    ,@(\"831\u00d72=1662\", \"478\u00d75=2390\")\n    ,@(\"199\u00d74=796\", \"840\u00d76=5040\")\n    ,@(\"849\u00d74=3396\", \"575\u00d73=1725\")\n    ,@(\"795\u00d75=3975\", \"619\u00d79=5571\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
